$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "2025-04-28 11:30:33"
$ws.Range("B10").Value = 218
